$wb = $excel.ActiveWorkbook

# --- YDS sheet: season running-total play logs (append Week 16 values) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "2 1 1 2 6 2 2 0 2 2 1 1 1 9 2 12 8 13 3 6 1 3 7 3 3 -1 17 5 1 4 4 6 10 -3 4 -2 10 17 3 0 0 4 4 14 11 2 4 1 4 2 2 2 13 8 5 4 2 -7 3 2 1 4 1 2 -2 9 6 2 5 -1 5 -2 -2 2 2 6 13 0 2 8 5 2 -3 5 3 -3 2 7 3 5 4 9 3 7 8 0 6 1 1 1 2 1 3 5 5 0 5 1 3 6 3 1 6 -3 19 4 3 0 6 8 6 2 5 9 -1 1 8 4 4 3 2 5 5 18 2 -1 1 4 2 7 1 4 5 -2 5 3 16 3 2 0 5 3 4 4 4 3 1 11 9 4 5 5 15 6 5 4 -4 5 -2 -2 5 8 5 13 0 0 0 5 1 -1 7 5 39 -1 5 4 15 1 1 4 2 5 0 2 5 3 4 1 1 19 2 2 -2 8 -1 6 12 2 2 3 1 8 6 11 5 4 1 9 4 2 24 13 2 0 2 5 -3 1 1 1 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17 3 2 18 3 1 2 -1 7 2 1 1 14 3 7 2 8 5 6 8 11 1 2 2 3 1 5 3 5 7 4 1 -5 11 0 4 6 6 20 -2 2 1 8 3 3 52 1 6 3 14 3 3 1 4 5 0 14 9 4 12 -1 3 20 7 5 10 0 11 2 0 3 1 38 3 3 2"
$ws.Range("C2").Value = "7 35 10 14 10 8 -3 1 11 20 7 8 22 6 15 2 40 12 8 27 7 11 9 8 5 3 17 4 5 3 27 5 8 6 6 12 12 11 4 8 7 6 5 6 -5 14 5 3 22 7 4 10 23 4 13 13 9 -4 13 30 8 6 3 6 15 6 54 5 29 3 53 4 0 9 29 11 8 11 14 1 13 3 6 14 4 5 -6 7 23 27 8 15 11 6 10 20 12 13 7 -2 6 15 3 7 22 5 17 21 8 7 21 8 8 12 15 6 7 3 4 6 6 16 9 30 8 7 6 3 15 3 6 7 12 6 20 6 5 2 26 11 12 15 9 4 8 13 14 15 9 8 9 10 24 14 23 12 6 8 8 19 13 12 0 11 12 15 26 19 -5 18 26 7 19 20 10 19 -1 18 12 12 0 13 -3 14 12 5 22 12 10 14 6 8 4 11 4 6 7 3 3 2 9 10 4 21 18 14 28 19 14 12 11 11 6 9 -7 21 18 9 2 5 5 29 14 20 9 5 22 5 4 9 23 5 5 2 62 12 8 11 8 15 11 5 11 6 1 10 22 16 -3 46 3 7 -3 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11 8 3 2 22 16 11 9 1 6 6 27 3 17 14 14 8 13 16 2 14 19 4 11 6 23 1 9 24 10 4 19 2 4 4 6 4 22 2 7 12 5"
$ws.Range("B3").Value = "7 8 4 1 7 2 3 1 5 8 4 3 2 2 9 0 14 3 -6 7 15 0 18 4 1 -1 9 2 5 2 7 7 2 4 2 26 9 16 1 -1 1 5 -2 12 3 -1 -5 7 3 8 4 11 0 1 14 2 2 0 13 4 13 1 10 0 -1 3 2 -3 10 4 7 2 -1 1 -2 1 0 -2 9 2 2 6 -1 7 4 0 2 2 2 4 16 1 7 2 0 4 3 5 2 3 1 2 11 2 3 18 9 11 22 1 2 2 7 5 5 1 0 -4 10 1 6 1 3 9 0 3 1 4 0 8 4 3 17 0 5 2 1 4 1 1 3 -2 3 6 3 12 5 9 2 32 6 1 3 10 10 0 13 3 0 -2 15 5 2 4 9 6 3 -2 1 6 -1 -2 0 1 -1 -1 -1 1 6 8 0 -4 6 7 4 9 4 5 -1 9 4 34 2 3 12 3 13 21 1 3 6 0 1 9 11 2 5 2 0 3 25 10 78 1 0 0 15 0 -4 4 0 4 11 7 13 24 8 15 0 1 7 6 3 9 2 14 3 -1 0 1 -3 5 0 4 1 2 7 3 2 20 9 9 0 -3 6 4 0 2 2 4 3 4 6 3 3 6 0 5 1 4 7 1 1 3 7 2 7 4 3 1 3 4 30 2 -5 2 5 0 7 -1 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0 5 11 7 1 1 3 10 -1 2 4 7 1 3 6 0 8 9 2 16 5 1 0 -2 7 13 -2 2 2 -1 9 3 0 9 3 8 3 1 -1 2 -1 0 0 44 3 4 0 0 4 10 1 1 12 4 5 0 8 2 5 4 1 6 6 9 3 7 2 0 11 3 2 0 3 24 8 30 4 12 1 4 3 0 9 3 1 1 6 6 6 0 8 0 2 4 3 2 4 7 1 6 3 9 7 2 3 6 1 0 26"
$ws.Range("C3").Value = "11 13 -5 8 1 9 26 5 27 15 11 57 14 4 5 11 22 17 1 7 12 7 3 -2 7 5 5 2 7 -1 2 28 19 3 24 -1 10 10 3 11 8 9 32 3 -1 1 6 10 8 5 28 6 17 6 9 27 6 21 -3 9 7 24 14 4 31 27 5 11 14 16 4 7 10 0 12 7 10 9 1 23 9 14 18 9 15 12 11 6 6 2 20 -6 5 9 12 0 3 16 -2 14 18 22 16 2 13 5 22 17 5 17 6 10 6 12 0 13 17 6 7 5 12 9 2 39 4 15 5 7 28 25 3 24 5 4 19 3 15 28 9 11 4 8 4 1 9 20 4 26 3 22 8 46 22 29 28 6 3 13 0 3 9 16 4 26 2 46 9 21 8 54 8 10 11 6 21 10 19 9 8 6 13 13 9 20 12 1 10 2 17 3 12 11 12 27 31 2 28 19 7 11 17 15 4 7 7 12 3 6 8 6 6 16 57 12 23 31 49 25 8 36 20 43 3 26 1 19 7 5 8 13 3 4 7 7 6 9 1 8 65 15 8 0 8 1 7 6 5 7 6 6 8 20 13 1 20 40 10 11 5 4 -4 4 7 2 6 5 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1 19 4 2 15 13 6 8 6 18 19 6 21 7 26 5 5 1 9 37 15 12 4 27 20 20 7 1 17 9 11 1 1 4 12 9 9 8 17 10 11 1 7 5 10 25 19 11 13 22 18 17 8 7 6 6 4"

# --- OFF sheet: Home (row 2) and Road (row 3) season totals ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 168
$ws.Range("D2").Value = 10
$ws.Range("F2").Value = 38
$ws.Range("G2").Value = 59
$ws.Range("J2").Value = 23
$ws.Range("L2").Value = 301
$ws.Range("M2").Value = 183
$ws.Range("O2").Value = 26
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 507
$ws.Range("C3").Value = 179
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 114
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 24
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 52
$ws.Range("N3").Value = 17

# --- DEF sheet: Home (row 2) and Road (row 3) season totals ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 226
$ws.Range("D2").Value = 19
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = 58
$ws.Range("J2").Value = 32
$ws.Range("L2").Value = 262
$ws.Range("M2").Value = 185
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 535
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 142
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 95
$ws.Range("G3").Value = 29
$ws.Range("H3").Value = 21
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 54
$ws.Range("N3").Value = 21

# --- ST sheet: summary counts + running logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 59
$ws.Range("D2").Value = 55
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 18
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 4
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = "65 49 62 70 56 61 64 62 66 67 50 55 64 66 62 55 64 62 64"
$ws.Range("B5").Value = "21 21 18 29 16 31 26 18 25 18 3 15 16 22 28 21 5 13 23"
$ws.Range("B6").Value = "0 38 18 21 23 30 24 26 16 15 27 65 27 22 23 25 19 25 35 24 27 25 33 18 29 33 29 20 20 31 29 79 28 32 30 23 21 23 42 20 2 34 31 29 19"
$ws.Range("D3").Value = "56 50 65 40 37 57 42 45 50 49 37 45 48 51 49 56 40 52 40 40 59 56 36 54 49 51 45 51 57 48 44 46 45 38 43 33 38 43 44 29 40 33 52 51 60 47 49 54 48 52 50 29 37 63 55"
$ws.Range("D4").Value = "4 9 0 0 0 12 15 0 0 0 0 5 0 11 16 22 0 0 0 6 0 15 0 4 17 13 0 0 10 0 0 26 0 0 0 0 0 0 7 0 0 0 11 12 13 0 0 18 0 0 0 0 0 5 4"
$ws.Range("D5").Value = "0 0 15 0 0 0 0 6 0 0 0 0 0 0 0 14 0 18 0 0 0 0 0 0 0 16 0 0 8 0 9 0 19 0 0 17 0 0 0 0 0 28 5 0 11 0 11 12 12"

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("E2").Value = 11
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 9

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 10
$ws.Range("D2").Value = 13
$ws.Range("D4").Value = 9

